# Updates the cryptos list table (Price and Volume(1h) columns; several
# coin rows also shift up by one as the source ranking changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.747.41"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3
$ws.Range("D3").Value = "1.799.47"
$ws.Range("E3").Value = "  -1.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4469"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.78%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3685"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07348"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8607"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.59%  "

# Row 12
$ws.Range("D12").Value = "1.797.16"
$ws.Range("E12").Value = "  -1.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.637"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.38"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07069"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.281"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008686"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.51%  "

# Row 19
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.86%  "

# Row 21
$ws.Range("D21").Value = "26.774.75"
$ws.Range("E21").Value = "  -1.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.160"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.978"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.23%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.169"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.71%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.205"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.09%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.06%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08783"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7417"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.47%  "

# Row 32
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.159"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.57%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.458"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.915"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.54%  "

# Row 35
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9997"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.084"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.15%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01961"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05198"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.81%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5290"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.37%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.832"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.24%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.977"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.26%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1685"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.34%  "

# Row 43
$ws.Range("B43").Value = "Decentraland"
$ws.Range("C43").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5092"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.39%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.447"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.40%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.68%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.969"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.23%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.82%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.678"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.65%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9993"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.15%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06292"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.19%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9168"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.51%  "
